# Adds a new "2022-Q4" quarter sheet (with fresh fund-holding data) right after
# the "总计" (total) sheet, and inserts the corresponding summary row at the top
# of the "总计" sheet's data table. All the other quarter sheets (2022-Q3 ..
# 2020-Q4) keep their own content unchanged; they are simply pushed one tab to
# the right as a natural side effect of inserting the new sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell while forcing a *text* number format first
# so numeric-looking strings ("6.89", "002207", ...) are not silently coerced
# into real numbers (which would also lose leading zeros on fund codes).
# ---------------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

function Set-NumberCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $val
}

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for 2022-Q4 and push the existing
#    quarters down. The leading index column (A) is a plain 0-based counter,
#    so after the insert every following row's index needs to be bumped by 1.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$lastRow = 9
$wsTotal.Rows.Item(2).Insert()

# Re-number the index column (A) for the rows that just shifted down
# (old row r -> new row r+1, old index value was (r-2), new index is (r-1)).
$stop = $lastRow + 1
for ($r = $stop; $r -ge 3; $r--) {
    $newIndex = $r - 2
    Set-NumberCell $wsTotal $r 1 $newIndex
}

# Fill the brand-new 2022-Q4 summary row.
Set-NumberCell $wsTotal 2 1 0
Set-TextCell   $wsTotal 2 2 "2022-Q4"
Set-NumberCell $wsTotal 2 3 4
Set-NumberCell $wsTotal 2 4 1.31

# ---------------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet right after "总计", modelled on the
#    "2022-Q3" sheet's layout (same header row / style), with the fresh
#    per-fund holdings for the new quarter.
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

$wsQ4 = $wb.Worksheets.Add($null, $wsTotal)
$wsQ4.Name = "2022-Q4"

# Copy header row (B1:H1) formatting + text from the 2022-Q3 sheet.
$wsQ3.Range("B1:H1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4104)
$wsQ3.Range("A1:A5").Copy()
$wsQ4.Range("A1:A5").PasteSpecial(-4104)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 2; $c -le 8; $c++) {
    Set-TextCell $wsQ4 1 $c $headers[$c - 2]
}

# Per-fund rows for 2022-Q4.
$rows = @(
    @{ idx = 0; code = "002207"; name = "前海开源金银珠宝主题精选混合C"; scale = "6.89"; pos = "90.36"; ratio = "7.75"; value = "0.5340"; rank = 7 },
    @{ idx = 1; code = "001302"; name = "前海开源金银珠宝主题精选混合A"; scale = "3.55"; pos = "90.36"; ratio = "7.75"; value = "0.2751"; rank = 7 },
    @{ idx = 2; code = "003304"; name = "前海开源沪港深核心资源灵活配置混合A"; scale = "3.30"; pos = "90.48"; ratio = "7.68"; value = "0.2534"; rank = 6 },
    @{ idx = 3; code = "003305"; name = "前海开源沪港深核心资源灵活配置混合C"; scale = "3.17"; pos = "90.48"; ratio = "7.68"; value = "0.2435"; rank = 6 }
)

$r = 2
foreach ($row in $rows) {
    Set-NumberCell $wsQ4 $r 1 $row.idx
    Set-TextCell   $wsQ4 $r 2 $row.code
    Set-TextCell   $wsQ4 $r 3 $row.name
    Set-TextCell   $wsQ4 $r 4 $row.scale
    Set-TextCell   $wsQ4 $r 5 $row.pos
    Set-TextCell   $wsQ4 $r 6 $row.ratio
    Set-TextCell   $wsQ4 $r 7 $row.value
    Set-NumberCell $wsQ4 $r 8 $row.rank
    $r++
}

# Restore the original active/selected tab (the last sheet, "2020-Q4"), since
# adding a new sheet otherwise leaves it focused instead.
$wsLast = $wb.Worksheets.Item("2020-Q4")
$wsLast.Select()
$wsLast.Range("A1").Select()
